$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (BUT 1 -> BUT 2, and fix double space)
$ws.Range("G2").Value = "BUT 2 INFORMATIQUE"

# SEMESTRE 1 in G3 stays the same

# Clear the year and jury date lines
$ws.Range("G4").ClearContents()
$ws.Range("G5").ClearContents()

# Update the active/selected cell to G3
$ws.Range("G3").Select()
